# This workbook is a "rodent" demo dataset: rows 2-6 hold five unique
# student records, and rows 7-11 / 12-16 are duplicates of that same
# template data (used by a data-cleanup/dedup exercise). The edit:
#   1) Gives each duplicated row its own unique numeric "Client Id" in
#      column C (previously every duplicate just re-used the shared
#      string for the original Id) so every row is distinguishable.
#   2) Fixes First Name (D) on a few rows that had mismatched names.
#   3) Straightens out the Overdue Disease / Overdue Agent / Imms Given /
#      Unique ID / Disease(s) / Imms History columns (M:R) on the rows
#      whose data had drifted out of sync with the rest of that record.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 7 (Whisker Elementary / Squeak template slot) ---
$ws.Range("C7").Value = 1009876548
$ws.Range("D7").Value = "Nibble"
$ws.Range("M7").Value = "Measles,"
$ws.Range("N7").Value = "MMR,"
$ws.Range("O7").Value = "Jul 10, 2014 - DTaP-IPV-Hib, Jul 10, 2014 - Pneu-C-13, Sep 15, 2014 - DTaP-IPV-Hib, Nov 20, 2014 - rota-unspecified, Mar 2, 2015 - MMR, Mar 2, 2015 - Men-C-C, Aug 7, 2015 - Var, Oct 1, 2015 - DTaP-IPV-Hib, May 19, 2024 - Tdap-IPV,"
$ws.Range("P7").Value = "CHEESE WHEEL ACADEMY-1009876544"
$ws.Range("Q7").Value = "Measles (MMR)"
$ws.Range("R7").Value = "[2014 JUL 10: DTaP-IPV-Hib, Pneu-C-13] [2014 SEP 15: DTaP-IPV-Hib] [2014 NOV 20: rota-unspecified] [2015 MAR 02: MMR, Men-C-C] [2015 AUG 07: Var] [2015 OCT 01: DTaP-IPV-Hib] [2024 MAY 19: Tdap-IPV]"

# --- Row 8 (Cheese Wheel Academy / Nibble template slot) ---
$ws.Range("C8").Value = 1009876549
$ws.Range("M8").Value = "Varicella,"
$ws.Range("N8").Value = "Var, HPV-9, Men-C-ACYW-135,"
$ws.Range("O8").Value = "Aug 20, 2013 - DTaP-IPV-Hib, Aug 20, 2013 - Pneu-C-13, Aug 20, 2013 - rota-unspecified, Nov 18, 2013 - DTaP-IPV-Hib, Nov 18, 2013 - Pneu-C-13, Jan 25, 2014 - DTaP-IPV-Hib, May 12, 2014 - MMR, May 12, 2014 - Men-C-C, Oct 3, 2014 - Var, Apr 14, 2024 - Tdap-IPV,"
$ws.Range("P8").Value = "WHISKER ELEMENTARY-1009876543"
$ws.Range("Q8").Value = "Varicella (Var)"
$ws.Range("R8").Value = "[2013 AUG 20: DTaP-IPV-Hib, Pneu-C-13, rota-unspecified] [2013 NOV 18: DTaP-IPV-Hib, Pneu-C-13] [2014 JAN 25: DTaP-IPV-Hib] [2014 MAY 12: MMR, Men-C-C] [2014 OCT 03: Var] [2024 APR 14: Tdap-IPV, MMR-Var]"

# --- Row 9 (Burrow Public School / Scurry template slot) ---
$ws.Range("C9").Value = 1009876550

# --- Row 10 (Tunnel Academy / Whiskers template slot) ---
$ws.Range("C10").Value = 1009876551
$ws.Range("D10").Value = "Chisel"

# --- Row 11 (Nutcracker Academy / Chisel template slot) ---
$ws.Range("C11").Value = 1009876552
$ws.Range("D11").Value = "Nibble"

# --- Row 12 (Whisker Elementary / Squeak template slot, 2nd dup) ---
$ws.Range("C12").Value = 1009876553
$ws.Range("D12").Value = "Chisel"
$ws.Range("M12").Value = "Hepatitis B,"
$ws.Range("N12").Value = "HB,"
$ws.Range("O12").Value = "Jan 5, 2013 - DTaP-IPV-Hib, Jan 5, 2013 - rota-unspecified, Mar 7, 2013 - Pneu-C-13, May 9, 2013 - DTaP-IPV-Hib, Jun 11, 2013 - MMR, Oct 23, 2013 - Men-C-C, Feb 2, 2014 - Var, May 6, 2014 - Pneu-C-13, Sep 12, 2014 - DTaP-IPV-Hib, May 1, 2024 - Tdap-IPV,"
$ws.Range("P12").Value = "BURROW PUBLIC SCHOOL-1009876545"
$ws.Range("Q12").Value = "Hepatitis B (HB)"
$ws.Range("R12").Value = "[2013 JAN 05: DTaP-IPV-Hib, rota-unspecified] [2013 MAR 07: Pneu-C-13] [2013 MAY 09: DTaP-IPV-Hib] [2013 JUN 11: MMR] [2013 OCT 23: Men-C-C] [2014 FEB 02: Var] [2014 MAY 06: Pneu-C-13] [2014 SEP 12: DTaP-IPV-Hib] [2024 MAY 01: Tdap-IPV]"

# --- Row 13 (Cheese Wheel Academy / Nibble template slot, 2nd dup) ---
$ws.Range("C13").Value = 1009876554

# --- Row 14 (Burrow Public School / Scurry template slot, 2nd dup) ---
$ws.Range("C14").Value = 1009876555
$ws.Range("M14").Value = "Varicella,"
$ws.Range("N14").Value = "Var, HPV-9, Men-C-ACYW-135,"
$ws.Range("O14").Value = "Aug 20, 2013 - DTaP-IPV-Hib, Aug 20, 2013 - Pneu-C-13, Aug 20, 2013 - rota-unspecified, Nov 18, 2013 - DTaP-IPV-Hib, Nov 18, 2013 - Pneu-C-13, Jan 25, 2014 - DTaP-IPV-Hib, May 12, 2014 - MMR, May 12, 2014 - Men-C-C, Oct 3, 2014 - Var, Apr 14, 2024 - Tdap-IPV,"
$ws.Range("P14").Value = "WHISKER ELEMENTARY-1009876543"
$ws.Range("Q14").Value = "Varicella (Var)"
$ws.Range("R14").Value = "[2013 AUG 20: DTaP-IPV-Hib, Pneu-C-13, rota-unspecified] [2013 NOV 18: DTaP-IPV-Hib, Pneu-C-13] [2014 JAN 25: DTaP-IPV-Hib] [2014 MAY 12: MMR, Men-C-C] [2014 OCT 03: Var] [2024 APR 14: Tdap-IPV, MMR-Var]"

# --- Row 15 (Tunnel Academy / Whiskers template slot, 2nd dup) ---
$ws.Range("C15").Value = 1009876556

# --- Row 16 (Nutcracker Academy / Chisel template slot, 2nd dup) ---
$ws.Range("C16").Value = 1009876557

# Column C now holds plain numeric Client Ids, so widen it to fit.
$ws.Columns("C").ColumnWidth = 9.92

# Leave the selection where the author left off editing.
$ws.Activate()
$ws.Range("M14:R14").Select()
